$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I33").Value = 351.18182
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 351.18182
$ws.Range("L33").Value = 90
$ws.Range("M33").Value = -122.18182
$ws.Range("N33").Value = -548

$ws.Range("H70").Value = 3666.889
$ws.Range("J70").Value = 1560
$ws.Range("L70").Value = 4680
$ws.Range("N70").Value = -5220

$ws.Range("H73").Value = 3666.889
$ws.Range("J73").Value = 1560
$ws.Range("L73").Value = 4680
$ws.Range("N73").Value = -6552

$ws.Range("H74").Value = 4793
$ws.Range("I74").Value = 4793
$ws.Range("K74").Value = 4793
$ws.Range("M74").Value = -3857

$ws.Range("H77").Value = 4793
$ws.Range("I77").Value = 4793
$ws.Range("K77").Value = 23965
$ws.Range("M77").Value = -19285

$ws.Range("H100").Value = 4371
$ws.Range("I100").Value = 5449.5
$ws.Range("K100").Value = 5449.5
$ws.Range("M100").Value = -4908.5

$ws.Range("H132").Value = 5996.1055
$ws.Range("I132").Value = 6129.4443
$ws.Range("J132").Value = 3596
$ws.Range("K132").Value = 18388.3329
$ws.Range("L132").Value = 10788
$ws.Range("M132").Value = -15858.3329
$ws.Range("N132").Value = -15848

$ws.Range("H141").Value = 9425
$ws.Range("I141").Value = 9425
$ws.Range("K141").Value = 28275
$ws.Range("M141").Value = -23095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1835
$ws.Range("I2").Value = 1838.5
$ws.Range("K2").Value = 1838.5
$ws.Range("M2").Value = -1725.5

$ws.Range("H102").Value = 1673.9166
$ws.Range("I102").Value = 1735.1818
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1735.1818
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -113.1818000000001
$ws.Range("N102").Value = -4244

$ws.Range("H116").Value = 1835
$ws.Range("I116").Value = 1838.5
$ws.Range("K116").Value = 1838.5
$ws.Range("M116").Value = 455.5

$ws.Range("H131").Value = 85000
$ws.Range("J131").Value = 85000
$ws.Range("L131").Value = 85000
$ws.Range("N131").Value = -95080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1835
$ws.Range("I3").Value = 1838.5
$ws.Range("K3").Value = 1838.5
$ws.Range("M3").Value = -1724.5

$ws.Range("H20").Value = 2233.45
$ws.Range("I20").Value = 1566.1538
$ws.Range("J20").Value = 3472.7144
$ws.Range("K20").Value = 1566.1538
$ws.Range("L20").Value = 3472.7144
$ws.Range("M20").Value = -1319.1538
$ws.Range("N20").Value = -3966.7144

$ws.Range("H64").Value = 2775.75
$ws.Range("J64").Value = 3970.5715
$ws.Range("L64").Value = 3970.5715
$ws.Range("N64").Value = -4420.5715

$ws.Range("H67").Value = 2775.75
$ws.Range("J67").Value = 3970.5715
$ws.Range("L67").Value = 3970.5715
$ws.Range("N67").Value = -5530.5715

$ws.Range("H105").Value = 3433.8572
$ws.Range("I105").Value = 2643.6316
$ws.Range("K105").Value = 2643.6316
$ws.Range("M105").Value = -896.6316000000002

$ws.Range("H134").Value = 4744.4443
$ws.Range("I134").Value = 4744.4443
$ws.Range("K134").Value = 14233.3329
$ws.Range("M134").Value = -11698.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 29964
$ws.Range("J20").Value = 29964
$ws.Range("L20").Value = 29964
$ws.Range("N20").Value = -30436

$ws.Range("H30").Value = 29964
$ws.Range("J30").Value = 29964
$ws.Range("L30").Value = 29964
$ws.Range("N30").Value = -30146

$ws.Range("H31").Value = 2999
$ws.Range("I31").Value = 1998.5
$ws.Range("K31").Value = 1998.5
$ws.Range("M31").Value = -1703.5

$ws.Range("H34").Value = 2999
$ws.Range("I34").Value = 1998.5
$ws.Range("K34").Value = 1998.5
$ws.Range("M34").Value = -1796.5

$ws.Range("H107").Value = 1078.2858
$ws.Range("J107").Value = 1272.5
$ws.Range("L107").Value = 1272.5
$ws.Range("N107").Value = -5112.5

$ws.Range("H128").Value = 29964
$ws.Range("J128").Value = 29964
$ws.Range("L128").Value = 29964
$ws.Range("N128").Value = -39924

$ws.Range("H132").Value = 1792.8
$ws.Range("I132").Value = 1792.8
$ws.Range("K132").Value = 5378.4
$ws.Range("M132").Value = -2848.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1443510.6
$ws.Range("I4").Value = 15686.15
$ws.Range("K4").Value = 47058.45
$ws.Range("M4").Value = -46946.45

$ws.Range("H22").Value = 434.66666
$ws.Range("I22").Value = 151
$ws.Range("J22").Value = 1002
$ws.Range("K22").Value = 453
$ws.Range("L22").Value = 3006
$ws.Range("M22").Value = -284
$ws.Range("N22").Value = -3344

$ws.Range("H27").Value = 434.66666
$ws.Range("I27").Value = 151
$ws.Range("J27").Value = 1002
$ws.Range("K27").Value = 453
$ws.Range("L27").Value = 3006
$ws.Range("M27").Value = -351
$ws.Range("N27").Value = -3210

$ws.Range("H131").Value = 1933.5
$ws.Range("I131").Value = 1578
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 4734
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = 306
$ws.Range("N131").Value = -19080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 149.11111
$ws.Range("I2").Value = 163.14285
$ws.Range("K2").Value = 163.14285
$ws.Range("M2").Value = -50.14285000000001

$ws.Range("H4").Value = 996.5
$ws.Range("J4").Value = 990
$ws.Range("L4").Value = 990
$ws.Range("N4").Value = -1214

$ws.Range("H70").Value = 6086.75
$ws.Range("I70").Value = 6299.1665
$ws.Range("J70").Value = 5449.5
$ws.Range("K70").Value = 6299.1665
$ws.Range("L70").Value = 5449.5
$ws.Range("M70").Value = -6029.1665
$ws.Range("N70").Value = -5989.5

$ws.Range("H73").Value = 6086.75
$ws.Range("I73").Value = 6299.1665
$ws.Range("J73").Value = 5449.5
$ws.Range("K73").Value = 6299.1665
$ws.Range("L73").Value = 5449.5
$ws.Range("M73").Value = -5363.1665
$ws.Range("N73").Value = -7321.5

$ws.Range("H128").Value = 45197.4
$ws.Range("I128").Value = 35999
$ws.Range("J128").Value = 47497
$ws.Range("K128").Value = 35999
$ws.Range("L128").Value = 47497
$ws.Range("M128").Value = -31019
$ws.Range("N128").Value = -57457

$ws.Range("H136").Value = 103500
$ws.Range("I136").Value = 200000
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 600000
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -597450
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5449.8823
$ws.Range("I7").Value = 5479.375
$ws.Range("K7").Value = 5479.375
$ws.Range("M7").Value = -5367.375

$ws.Range("H22").Value = 1510
$ws.Range("I22").Value = 1440
$ws.Range("K22").Value = 1440
$ws.Range("M22").Value = -1145

$ws.Range("H27").Value = 1510
$ws.Range("I27").Value = 1440
$ws.Range("K27").Value = 1440
$ws.Range("M27").Value = -1333

$ws.Range("H68").Value = 2899.6
$ws.Range("J68").Value = 3349.5
$ws.Range("L68").Value = 3349.5
$ws.Range("N68").Value = -4847.5

$ws.Range("H71").Value = 2899.6
$ws.Range("J71").Value = 3349.5
$ws.Range("L71").Value = 16747.5
$ws.Range("N71").Value = -24235.5

$ws.Range("H82").Value = 794.44446
$ws.Range("J82").Value = 510.8
$ws.Range("L82").Value = 510.8
$ws.Range("N82").Value = -1232.8

$ws.Range("H85").Value = 794.44446
$ws.Range("J85").Value = 510.8
$ws.Range("L85").Value = 510.8
$ws.Range("N85").Value = -3006.8

$ws.Range("H126").Value = 5449.8823
$ws.Range("I126").Value = 5479.375
$ws.Range("K126").Value = 16438.125
$ws.Range("M126").Value = -13968.125

$ws.Range("H128").Value = 79282
$ws.Range("J128").Value = 79282
$ws.Range("L128").Value = 79282
$ws.Range("N128").Value = -89242

$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 11099.6
$ws.Range("J107").Value = 9374.5
$ws.Range("L107").Value = 28123.5
$ws.Range("N107").Value = -31963.5

$ws.Range("H130").Value = 55497.5
$ws.Range("J130").Value = 55497.5
$ws.Range("L130").Value = 55497.5
$ws.Range("N130").Value = -65537.5
